$wb = $excel.ActiveWorkbook

# --- Sheet: SoCDTtiNTY-psgr ---
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

$wsPsgr.Range("B2").Value = 0.36973314103415234
$wsPsgr.Range("C2").Value = 0.072
$wsPsgr.Range("D2").Value = 0.08755837753596755
$wsPsgr.Range("E2").Value = 0.08755837753596757
$wsPsgr.Range("F2").Value = 0.072
$wsPsgr.Range("G2").Value = 0.072
$wsPsgr.Range("H2").Value = 0.8235507771969761

$wsPsgr.Range("B3").Value = 0.9501394088209177
$wsPsgr.Range("C3").Value = 0.07464627318333034
$wsPsgr.Range("D3").Value = 0.07464627318333034
$wsPsgr.Range("E3").Value = 0.07464627318333034
$wsPsgr.Range("F3").Value = 0.07464627318333034
$wsPsgr.Range("G3").Value = 0.07464627318333034
$wsPsgr.Range("H3").Value = 0.07464627318333034

$wsPsgr.Range("B4").Value = 0.08
$wsPsgr.Range("C4").Value = 0.08
$wsPsgr.Range("D4").Value = 0.08
$wsPsgr.Range("E4").Value = 0.08
$wsPsgr.Range("F4").Value = 0.08
$wsPsgr.Range("G4").Value = 0.08
$wsPsgr.Range("H4").Value = 0.08

# --- Sheet: SoCDTtiNTY-frgt ---
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

$wsFrgt.Range("B2").Value = 0.36973314103415234
$wsFrgt.Range("D2").Value = 0.08755837753596755
$wsFrgt.Range("E2").Value = 0.08755837753596757
$wsFrgt.Range("H2").Value = 0.8235507771969761

$wsFrgt.Range("B3").Value = 0.9501394088209177
$wsFrgt.Range("C3").Value = 0.07464627318333034
$wsFrgt.Range("D3").Value = 0.07464627318333034
$wsFrgt.Range("E3").Value = 0.07464627318333034
$wsFrgt.Range("F3").Value = 0.07464627318333034
$wsFrgt.Range("G3").Value = 0.07464627318333034
$wsFrgt.Range("H3").Value = 0.07464627318333034

$wsFrgt.Range("B4").Value = 0.08
$wsFrgt.Range("C4").Value = 0.08
$wsFrgt.Range("D4").Value = 0.08
$wsFrgt.Range("E4").Value = 0.08
$wsFrgt.Range("F4").Value = 0.08
$wsFrgt.Range("G4").Value = 0.08
$wsFrgt.Range("H4").Value = 0.08

# --- Selection / active sheet adjustments ---
$wsFrgt.Activate()
$wsFrgt.Range("C14").Select()
$wsPsgr.Activate()
$wsPsgr.Range("C14").Select()
